$wb = $excel.ActiveWorkbook

$oldId = "9e7a7080-8f21-46a8-81bb-a6ea28049ead"
$newId = "85359eea-a0d9-469c-b494-81ad2217908b"
$oldHash = "eab11fba828a382dd17d3c8aabea4840849e434d"
$newHash = "6088bf3042ce6054d453aa35b926c0d4231dcf15"

$newMd = "$newId.md"
$newZh = "$newId.$newHash.zh-cn.xlf"
$newDe = "$newId.$newHash.de-de.xlf"

$newHandoffDate = "2016-03-24 01:05:59"
$newZhDatetime = "2016-03-24 01:05:55"

$mdAddress = "https://github.com/OpenLocalizationTest/oltest/blob/cfd647989e4bcb6c5fe34515089f32619e2b2ce3/e2e/$oldId.md"
$zhAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/996d791762b4e6208bd7279d1ffd1737688aa691/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldId.$oldHash.zh-cn.xlf"
$deAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/778e5939ba1bff6ee89db2d785f81c57a88bd9a1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldId.$oldHash.de-de.xlf"

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMd
$wsOverview.Range("D2").Value = $newHandoffDate

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $mdAddress, "", "", $newMd)

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $newMd
$wsZh.Range("D2").Value = $newZh
$wsZh.Range("E2").Value = $newZhDatetime

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $mdAddress, "", "", $newMd)
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), $zhAddress, "", "", $newZh)

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $newMd
$wsDe.Range("D2").Value = $newDe
$wsDe.Range("E2").Value = $newHandoffDate

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $mdAddress, "", "", $newMd)
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), $deAddress, "", "", $newDe)
